$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yellow = 65535

# --- Column A: renamed / new rows (11 and 12) ---
$ws.Range("A11").Value = "Consultar mis datos"
$ws.Range("A11").Interior.Color = $yellow

$ws.Range("A12").Value = "Incrementar Cantidad de Eventos"

# --- Column B: rename "Borrar" rows to "Actuvar o Desactivar", add new rows ---
$ws.Range("B9").Value = "Actuvar o Desactivar Candidato"
$ws.Range("B10").Value = "Actuvar o Desactivar Participante"

$ws.Range("B14").Value = "Registro Masivo"
$ws.Range("B14").Interior.Color = $yellow
$ws.Range("B15").Value = "Carga imagen del Evento"
$ws.Range("B15").Interior.Color = $yellow

# --- Column C: add "Registrame en eventos por código" row + blank styled cells ---
$ws.Range("C8").Value = "Registrame en eventos por código"
$ws.Range("C8").Interior.Color = $yellow

$ws.Range("C11").Interior.Color = $yellow
$ws.Range("C12").Interior.Color = $yellow
$ws.Range("C13").Interior.Color = $yellow
$ws.Range("C14").Interior.Color = $yellow

# --- Column D: rename "Borrar" rows to "Actuvar o Desactivar", add new row ---
$ws.Range("D7").Value = "Actuvar o Desactivar Administradores"
$ws.Range("D8").Value = "Actuvar o Desactivar Evento"
$ws.Range("D9").Value = "Actuvar o Desactivar Candidato"
$ws.Range("D10").Value = "Actuvar o Desactivar Participante"

$ws.Range("D14").Value = "Registro Masivo"
$ws.Range("D14").Interior.Color = $yellow

# --- Column E: add image-related rows ---
$ws.Range("E8").Value = "Subir Imágenes"
$ws.Range("E8").Interior.Color = $yellow
$ws.Range("E9").Value = "Borrar Imágenes"
$ws.Range("E9").Interior.Color = $yellow

# --- Fill-color-only additions on existing cells (no value change) ---
$ws.Range("B6:E6").Interior.Color = $yellow
$ws.Range("C7").Interior.Color = $yellow
$ws.Range("E7").Interior.Color = $yellow
$ws.Range("B11").Interior.Color = $yellow
$ws.Range("D11").Interior.Color = $yellow
$ws.Range("B12").Interior.Color = $yellow
$ws.Range("D12").Interior.Color = $yellow
$ws.Range("B13").Interior.Color = $yellow
$ws.Range("D13").Interior.Color = $yellow

# --- Selection matches the saved view state ---
$ws.Range("A13").Select()
